$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "# - Mutual Friends count"

$ws.Activate() | Out-Null
$ws.Range("B1:G1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 145
